# Matriz EDT.xlsx — apply the "avance matriz EDT" edit:
#  - EDT!B5 gets the "dias " column sub-header label.
#  - EDT!B6:B56 (excluding the phase-header rows 11/18/47/52) become
#    formulas =(C+D+E+F+G)/8, formatted as "0.0".
#  - The phase-header rows (11/18/47/52) keep their B cell blank but pick
#    up the same "0.0" numeric style (and lose their grey banding fill,
#    matching the source edit).
#  - View/selection state: EDT scrolled to A29 with I17 selected and kept
#    as the active sheet; "Costos Fases" selection parked at E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDT")
$ws2 = $wb.Worksheets.Item("Costos Fases")

# --- EDT!B5: new "dias " sub-header label ------------------------------
$ws.Range("B5").Value = "dias "

# --- EDT!B6:B56: per-row hours/8 formula --------------------------------
$headerRows = @(11, 18, 47, 52)
for ($r = 6; $r -le 56; $r++) {
    if ($headerRows -contains $r) { continue }
    $formula = "=(C" + $r + "+D" + $r + "+E" + $r + "+F" + $r + "+G" + $r + ")/8"
    $ws.Cells.Item($r, 2).Formula = $formula
}

# --- Number format "0.0" across the whole DIAS column (B6:B56) ---------
$ws.Range("B6:B56").NumberFormat = "0.0"

# The phase-header rows (B11/B18/B47/B52) drop their grey fill once the
# new numeric style is applied, matching the authored edit.
foreach ($r in $headerRows) {
    $ws.Cells.Item($r, 2).Interior.Pattern = -4142   # xlNone
}

# --- Recalculate so B58's SUM(B6:B57) reflects the new formulas --------
$excel.Calculate()

# --- View / selection state ---------------------------------------------
# Touch "Costos Fases" selection first so it doesn't end up as the active
# tab once we're done (EDT must stay tabSelected).
$ws2.Select()
$ws2.Range("E3").Select()

$ws.Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
$ws.Range("I17").Select()
